{"js": "// Add new notes on tensors, gradients and derivatives to the end of the\n// document, after the existing \"Now form a set...\" paragraph, and move the\n// trailing \"_GoBack\" bookmark so it still marks the end of the document.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst lastParagraph = paragraphs.items[paragraphs.items.length - 1];\n\n// The \"_GoBack\" bookmark currently sits between the trailing space run and\n// the \"Now form a set...\" run in the last paragraph; remove it from there -\n// it will be re-inserted at the end of the new content below.\ncontext.document.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\n// New paragraphs to append, expressed as raw OOXML so run/formatting\n// details (tab, superscript, spell-check markers) come through exactly.\nconst newParagraphsOoxml =\n  \"<w:p/>\" +\n  \"<w:p><w:r><w:t>What is tensor?</w:t></w:r></w:p>\" +\n  \"<w:p><w:r><w:tab/><w:t>Tensor can be anything, from a number, to an array, to a vector, to a vector of vector of vectors aka n-dimensional array.</w:t></w:r></w:p>\" +\n  \"<w:p><w:r><w:t>Matrix multiplication along with other operations such as addition, division, etc. can be done easily using tensors.</w:t></w:r></w:p>\" +\n  \"<w:p>\" +\n    \"<w:r><w:t xml:space=\\\"preserve\\\">Tensor maintains conformity. </w:t></w:r>\" +\n    \"<w:r><w:t>If in the 4</w:t></w:r>\" +\n    \"<w:r><w:rPr><w:vertAlign w:val=\\\"superscript\\\"/></w:rPr><w:t>th</w:t></w:r>\" +\n    \"<w:r><w:t xml:space=\\\"preserve\\\"> dimension a row has 12 columns, all the rows will have equal number of columns. If one of those columns have 7 rows than a</w:t></w:r>\" +\n    \"<w:r><w:t>ll 12 columns will have 7 rows.</w:t></w:r>\" +\n  \"</w:p>\" +\n  \"<w:p>\" +\n    \"<w:r><w:t xml:space=\\\"preserve\\\">The term gradients and derivatives are same. It mathematical expression is </w:t></w:r>\" +\n    \"<w:proofErr w:type=\\\"spellStart\\\"/>\" +\n    \"<w:r><w:t>dy</w:t></w:r>\" +\n    \"<w:proofErr w:type=\\\"spellEnd\\\"/>\" +\n    \"<w:r><w:t xml:space=\\\"preserve\\\">/dx and it is the slope of a curve. </w:t></w:r>\" +\n    \"<w:r><w:t xml:space=\\\"preserve\\\">Gradient or derivative is for 1 unit of movement along with x-axis, how much movement happens in y-axis. </w:t></w:r>\" +\n  \"</w:p>\" +\n  \"<w:p><w:r><w:t xml:space=\\\"preserve\\\">Derivatives are used when you are dealing with numbers and gradients are used when you are dealing with matrices. </w:t></w:r></w:p>\";\n\nconst packageOoxml =\n  \"<pkg:package xmlns:pkg=\\\"http://schemas.microsoft.com/office/2006/xmlPackage\\\">\" +\n    \"<pkg:part pkg:name=\\\"/word/document.xml\\\" pkg:contentType=\\\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\\\">\" +\n      \"<pkg:xmlData>\" +\n        \"<w:document xmlns:w=\\\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\\\">\" +\n          \"<w:body>\" + newParagraphsOoxml + \"</w:body>\" +\n        \"</w:document>\" +\n      \"</pkg:xmlData>\" +\n    \"</pkg:part>\" +\n  \"</pkg:package>\";\n\nconst insertionPoint = lastParagraph.getRange(Word.RangeLocation.end);\ninsertionPoint.insertOoxml(packageOoxml, Word.InsertLocation.after);\nawait context.sync();\n\n// Put the \"_GoBack\" bookmark back at the very end of the document (end of\n// the last, newly-inserted paragraph).\nconst refreshedParagraphs = body.paragraphs;\nrefreshedParagraphs.load(\"items\");\nawait context.sync();\n\nconst newLastParagraph = refreshedParagraphs.items[refreshedParagraphs.items.length - 1];\nconst endOfDoc = newLastParagraph.getRange(Word.RangeLocation.end);\nendOfDoc.insertBookmark(\"_GoBack\");\nawait context.sync();\n", "ps1": "# Add new notes on tensors, gradients and derivatives to the end of the\n# document, after the existing \"Now form a set...\" paragraph, and move the\n# trailing \"_GoBack\" bookmark so it still marks the end of the document.\n\n$d = $word.ActiveDocument\n\n# The \"_GoBack\" bookmark currently sits between the trailing space run and\n# the \"Now form a set...\" run in the last paragraph; remove it from there -\n# it will be re-inserted (via the OOXML below) at the end of the new\n# content.\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks(\"_GoBack\").Delete()\n}\n\n$lastParagraph = $d.Paragraphs.Last\n\n# New paragraphs to append, expressed as raw OOXML so run/formatting\n# details (tab, superscript, spell-check markers, bookmark) come through\n# exactly. The trailing bookmarkStart/bookmarkEnd re-establish \"_GoBack\" at\n# the new end of the document.\n$w = \"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"\n$newParagraphsOoxml = (\n    \"<w:p xmlns:w=`\"$w`\"/>\" +\n    \"<w:p xmlns:w=`\"$w`\"><w:r><w:t>What is tensor?</w:t></w:r></w:p>\" +\n    \"<w:p xmlns:w=`\"$w`\"><w:r><w:tab/><w:t>Tensor can be anything, from a number, to an array, to a vector, to a vector of vector of vectors aka n-dimensional array.</w:t></w:r></w:p>\" +\n    \"<w:p xmlns:w=`\"$w`\"><w:r><w:t>Matrix multiplication along with other operations such as addition, division, etc. can be done easily using tensors.</w:t></w:r></w:p>\" +\n    \"<w:p xmlns:w=`\"$w`\">\" +\n        \"<w:r><w:t xml:space=`\"preserve`\">Tensor maintains conformity. </w:t></w:r>\" +\n        \"<w:r><w:t>If in the 4</w:t></w:r>\" +\n        \"<w:r><w:rPr><w:vertAlign w:val=`\"superscript`\"/></w:rPr><w:t>th</w:t></w:r>\" +\n        \"<w:r><w:t xml:space=`\"preserve`\"> dimension a row has 12 columns, all the rows will have equal number of columns. If one of those columns have 7 rows than a</w:t></w:r>\" +\n        \"<w:r><w:t>ll 12 columns will have 7 rows.</w:t></w:r>\" +\n    \"</w:p>\" +\n    \"<w:p xmlns:w=`\"$w`\">\" +\n        \"<w:r><w:t xml:space=`\"preserve`\">The term gradients and derivatives are same. It mathematical expression is </w:t></w:r>\" +\n        \"<w:proofErr w:type=`\"spellStart`\"/>\" +\n        \"<w:r><w:t>dy</w:t></w:r>\" +\n        \"<w:proofErr w:type=`\"spellEnd`\"/>\" +\n        \"<w:r><w:t xml:space=`\"preserve`\">/dx and it is the slope of a curve. </w:t></w:r>\" +\n        \"<w:r><w:t xml:space=`\"preserve`\">Gradient or derivative is for 1 unit of movement along with x-axis, how much movement happens in y-axis. </w:t></w:r>\" +\n    \"</w:p>\" +\n    \"<w:p xmlns:w=`\"$w`\"><w:r><w:t xml:space=`\"preserve`\">Derivatives are used when you are dealing with numbers and gradients are used when you are dealing with matrices. </w:t></w:r><w:bookmarkStart w:id=`\"0`\" w:name=`\"_GoBack`\"/><w:bookmarkEnd w:id=`\"0`\"/></w:p>\"\n)\n\n# Collapsed, zero-length range positioned right before the paragraph mark\n# that ends the last paragraph, so the new paragraphs are inserted after it\n# without consuming/merging that paragraph mark.\n$insertionPoint = $d.Range($lastParagraph.Range.End - 1, $lastParagraph.Range.End - 1)\n$null = $insertionPoint.InsertXML($newParagraphsOoxml)\n"}
